$d = $word.ActiveDocument

$replacements = @(
    @{old="78×51=3978"; new="19×35=665"},
    @{old="85×67=5695"; new="87×25=2175"},
    @{old="36×38=1368"; new="75×83=6225"},
    @{old="40×42=1680"; new="39×64=2496"},
    @{old="14×80=1120"; new="12×86=1032"},
    @{old="16×33=528";  new="36×85=3060"},
    @{old="20×98=1960"; new="38×16=608"},
    @{old="17×77=1309"; new="59×18=1062"},
    @{old="81×72=5832"; new="43×14=602"},
    @{old="54×15=810";  new="55×18=990"},
    @{old="30×57=1710"; new="52×72=3744"},
    @{old="96×78=7488"; new="82×40=3280"},
    @{old="55×78=4290"; new="86×31=2666"},
    @{old="97×27=2619"; new="48×34=1632"},
    @{old="49×75=3675"; new="11×34=374"},
    @{old="17×75=1275"; new="69×57=3933"},
    @{old="95×28=2660"; new="85×57=4845"},
    @{old="75×13=975";  new="99×13=1287"},
    @{old="30×43=1290"; new="28×71=1988"},
    @{old="57×32=1824"; new="72×71=5112"},
    @{old="50×90=4500"; new="75×78=5850"},
    @{old="72×26=1872"; new="70×69=4830"},
    @{old="73×37=2701"; new="54×65=3510"},
    @{old="46×82=3772"; new="66×65=4290"},
    @{old="96×85=8160"; new="87×92=8004"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
